$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 319 (shifts the existing 319-337 data down to 321-339)
$ws.Rows("319:320").Insert()

# New row 319: Angeleno / Primera, week of 2023-04-25 (serial 45041)
$ws.Range("A319").Value = 4
$ws.Range("B319").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C319").Value = "Los Lagos"
$ws.Range("D319").Value = 45041
$ws.Range("E319").Value = 10
$ws.Range("F319").Value = "Fruta"
$ws.Range("G319").Value = 100103
$ws.Range("H319").Value = "Frutos de hueso (carozo)"
$ws.Range("I319").Value = 100103002
$ws.Range("J319").Value = "Ciruela"
$ws.Range("K319").Value = "Angeleno"
$ws.Range("L319").Value = "Primera"
$ws.Range("M319").Value = 500
$ws.Range("N319").Value = 15000
$ws.Range("O319").Value = 16000
$ws.Range("P319").Value = 15600
$ws.Range("Q319").Value = "$/caja 14 kilos granel"
$ws.Range("R319").Value = "Región de O'Higgins"
$ws.Range("S319").Value = 1114
$ws.Range("T319").Value = 14

# New row 320: Angeleno / Segunda, same week
$ws.Range("A320").Value = 4
$ws.Range("B320").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C320").Value = "Los Lagos"
$ws.Range("D320").Value = 45041
$ws.Range("E320").Value = 10
$ws.Range("F320").Value = "Fruta"
$ws.Range("G320").Value = 100103
$ws.Range("H320").Value = "Frutos de hueso (carozo)"
$ws.Range("I320").Value = 100103002
$ws.Range("J320").Value = "Ciruela"
$ws.Range("K320").Value = "Angeleno"
$ws.Range("L320").Value = "Segunda"
$ws.Range("M320").Value = 300
$ws.Range("N320").Value = 14000
$ws.Range("O320").Value = 14000
$ws.Range("P320").Value = 14000
$ws.Range("Q320").Value = "$/caja 14 kilos granel"
$ws.Range("R320").Value = "Región de O'Higgins"
$ws.Range("S320").Value = 1000
$ws.Range("T320").Value = 14
